$d = $word.ActiveDocument

# Title / H1 heading AND the bold "recap" run near the end both share this
# exact text, so a single ReplaceAll pass over the whole document updates
# both occurrences in one shot.
$d.Content.Find.Execute(
    "Play Book of Itza Free - Innovative Aztec Themed Slot", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Play Book of Itza Free: Innovative Gameplay & Exciting Bonus Features", 2)

# "What we like" bullet list
$d.Content.Find.Execute(
    "Innovative gameplay with 1,024 ways to win", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Innovative mechanism for forming combinations", 2)

$d.Content.Find.Execute(
    "Aztec theme with well-designed graphics and pleasing soundtrack", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Well-designed graphics and clean visuals", 2)

$d.Content.Find.Execute(
    "Special symbols and lucrative bonus features", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Exciting bonus features like Expanded Wild Respins and free spins", 2)

$d.Content.Find.Execute(
    "Medium volatility and high RTP for a fair chance to win", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Balanced gameplay with medium volatility", 2)

# "What we don't like" bullet list
$d.Content.Find.Execute(
    "Expanding Wild is the highest point of luck", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Limited variety in terms of symbols and theme", 2)

$d.Content.Find.Execute(
    "Free spin feature is difficult to trigger", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Medium volatility may not appeal to high-risk players", 2)

# Meta description (italic run near the end)
$d.Content.Find.Execute(
    "Discover a balanced and enjoyable gaming experience with innovative gameplay mechanics and lucrative bonus features in Book of Itza online slot game. Play for free.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Discover Book of Itza, an online slot game with innovative gameplay and exciting bonus features. Play for free now!", 2)
